$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59's phone number had been recorded as text; it now gets re-entered
# as a plain number.
$ws.Cells.Item(59, 1).Value = 76442781

# Append the new payment (76442781, Cash, 2025-08-18T18:07:47) as row 60.
# The phone number is text (a leading apostrophe forces text so the
# numeric-looking string isn't auto-converted to a number); the blank
# method/discount columns use a lone apostrophe so they stay empty text
# cells instead of disappearing entirely. The style is reset back to
# Normal afterwards so the quote-prefix marker doesn't linger as a style.
$ws.Cells.Item(60, 1).Value = "'76442781"
$ws.Cells.Item(60, 1).Style = "Normal"
$ws.Cells.Item(60, 2).Value = "'"
$ws.Cells.Item(60, 2).Style = "Normal"
$ws.Cells.Item(60, 3).Value = "Cash"
$ws.Cells.Item(60, 4).Value = "2025-08-18T18:07:47"
$ws.Cells.Item(60, 5).Value = 120
$ws.Cells.Item(60, 6).Value = "'"
$ws.Cells.Item(60, 6).Style = "Normal"
$ws.Cells.Item(60, 7).Value = 120
$ws.Cells.Item(60, 8).Value = 0
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 0
